# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G (header "K") values for rows 2-29 with the newly computed
# strikeout counts (replacing the previous "Strike#" values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 2
    6  = 1
    7  = 1
    8  = 1
    9  = 0
    10 = 1
    11 = 2
    12 = 0
    13 = 1
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 1
    19 = 0
    20 = 1
    21 = 0
    22 = 2
    23 = 1
    24 = 2
    25 = 1
    26 = 0
    27 = 2
    28 = 1
    29 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
